$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "S2"  = 0.2118068965517241
    "S3"  = 0.1805458229957766
    "S7"  = 0.2757242757242757
    "S10" = 0.216893039049236
    "S11" = 0.2738805263656158
    "S13" = 0.2939434534301629
    "S14" = 0.2832591683289857
    "S15" = 0.2672127950068266
    "S16" = 0.2971181376820615
    "S19" = 0.2493173885859726
    "S24" = 0.2075536661425302
    "S25" = 0.2711942992381789
    "S26" = 0.2733596891596734
    "S27" = 0.2498671219559148
    "S29" = 0.2504936941240743
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
